# Timesheet update: fix/add clock-in & clock-out entries for early Feb 2026.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: recalculated duration ---
$ws.Range("D10").Value = "1.08 Hours"

# --- Row 12: fill in the clock-out time + computed duration ---
$ws.Range("C12").Value = "23:03:37"
$ws.Range("D12").Value = "1.93 Hours"

# --- Row 13: brand-new clock-in entry (clock-out not punched yet) ---
# Copy the style (font/format) of an existing data row onto row 13 first,
# so the new cells line up with the rest of the table (s="2") instead of
# picking up default formatting / auto-detected number formats.
$ws.Range("A10:D10").Copy()
$ws.Range("A13:D13").PasteSpecial(-4122)

# A13 holds a literal date-like string ("2026-02-03"). Entering that via
# .Value would be auto-parsed into a real date serial (like typing it into
# Excel would). Build the literal text on a scratch cell via a formula
# (so it is a plain text result, not a date), then paste just the VALUE
# into A13 - this keeps the pasted-in style/format untouched.
$ws.Range("Z1").Formula = "=""2026-02-03"""
$ws.Range("Z1").Copy()
$ws.Range("A13").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

# B13 is a plain time-of-day string; it is not auto-converted, so a
# straightforward value assignment is fine.
$ws.Range("B13").Value = "07:59:16"

# C13 / D13 stay blank (no clock-out yet), matching C12/D12's prior state -
# already handled by the format-only paste above (blank numeric cell, s="2").

# --- Row heights: rows 10-13 become 15.75pt custom height rows ---
$ws.Range("A10:D13").RowHeight = 15.75

# --- Selection moves to D11 ---
$ws.Range("D11").Select()
